$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in grades (value 5) for row 11 (columns C, D, E)
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 5
$ws.Range("E11").Value = 5

# Fill in grades (value 5) for row 13 (columns C, D, E, F)
$ws.Range("C13").Value = 5
$ws.Range("D13").Value = 5
$ws.Range("E13").Value = 5
$ws.Range("F13").Value = 5

# Update the active selection / cell to G13
$ws.Range("G13").Select()
